$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 55. This shifts existing rows 55-107 down to 56-108,
# carrying along their values/styles, matching the dimension change to A1:R108.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with a new weekly price record.
$ws.Cells.Item(55, 1).Value = 10
$ws.Cells.Item(55, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(55, 3).Value = "La Araucanía"
$ws.Cells.Item(55, 4).Value = 44778
$ws.Cells.Item(55, 4).Style = $ws.Cells.Item(56, 4).Style
$ws.Cells.Item(55, 4).NumberFormat = $ws.Cells.Item(56, 4).NumberFormat
$ws.Cells.Item(55, 5).Value = 9
$ws.Cells.Item(55, 6).Value = 100112035
$ws.Cells.Item(55, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 85
$ws.Cells.Item(55, 11).Value = 25000
$ws.Cells.Item(55, 12).Value = 25000
$ws.Cells.Item(55, 13).Value = 25000
$ws.Cells.Item(55, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(55, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(55, 16).Value = 2500
$ws.Cells.Item(55, 17).Value = 10
$ws.Cells.Item(55, 18).Value = "Hortaliza"
